# The workbook tracks daily Brócoli prices for "Femacal de La Calera".
# This edit adds one new reporting day (two quality rows: Primera / Segunda)
# at its chronological slot (rows 389-390), pushing the existing rows
# (389-481) down by two and growing the used range to A1:R483.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 389, shifting rows 389:481 down to 391:483.
$ws.Rows("389:390").Insert()

# Populate the two newly inserted rows with the new day's data.
# Row 389 - Primera
$ws.Cells.Item(389, 1).Value  = 3
$ws.Cells.Item(389, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(389, 3).Value  = "Coquimbo"
$ws.Cells.Item(389, 4).Value  = 44543
$ws.Cells.Item(389, 5).Value  = 5
$ws.Cells.Item(389, 6).Value  = 100112023
$ws.Cells.Item(389, 7).Value  = "Brócoli"
$ws.Cells.Item(389, 8).Value  = "Sin especificar"
$ws.Cells.Item(389, 9).Value  = "Primera"
$ws.Cells.Item(389, 10).Value = 1300
$ws.Cells.Item(389, 11).Value = 600
$ws.Cells.Item(389, 12).Value = 600
$ws.Cells.Item(389, 13).Value = 600
$ws.Cells.Item(389, 14).Value = "`$/unidad"
$ws.Cells.Item(389, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(389, 16).Value = 600
$ws.Cells.Item(389, 17).Value = 1
$ws.Cells.Item(389, 18).Value = "Hortaliza"

# Row 390 - Segunda
$ws.Cells.Item(390, 1).Value  = 3
$ws.Cells.Item(390, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(390, 3).Value  = "Coquimbo"
$ws.Cells.Item(390, 4).Value  = 44543
$ws.Cells.Item(390, 5).Value  = 5
$ws.Cells.Item(390, 6).Value  = 100112023
$ws.Cells.Item(390, 7).Value  = "Brócoli"
$ws.Cells.Item(390, 8).Value  = "Sin especificar"
$ws.Cells.Item(390, 9).Value  = "Segunda"
$ws.Cells.Item(390, 10).Value = 2700
$ws.Cells.Item(390, 11).Value = 400
$ws.Cells.Item(390, 12).Value = 500
$ws.Cells.Item(390, 13).Value = 452
$ws.Cells.Item(390, 14).Value = "`$/unidad"
$ws.Cells.Item(390, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(390, 16).Value = 452
$ws.Cells.Item(390, 17).Value = 1
$ws.Cells.Item(390, 18).Value = "Hortaliza"

# Make sure the date column keeps the same date-time number format as the
# rest of column D (style index picked up "s=2" automatically via Insert,
# but set format explicitly too in case it didn't carry over).
$ws.Range("D389:D390").NumberFormat = $ws.Range("D388").NumberFormat
